$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.851.47"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "1.857.71"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5037"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3624"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07162"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8919"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "1.840.63"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.218"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008492"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "26.887.50"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.012"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").Value = "2.076.28"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.405"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.792"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.044"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.643"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.653"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09229"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7438"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.975"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.147"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.254"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.505"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01985"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.088"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5336"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "118.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.493"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.437"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1462"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.0000"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05934"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.62%  "
